$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

$ws.Range("B25").Value = "acted"
$ws.Range("C25").Value = 201
$ws.Range("D25").Value = 1771
$ws.Range("E25").Value = 1028
$ws.Range("F25").Value = "Persian"
$ws.Range("G25").Formula = "=IF(OR(ISBLANK(C25), ISBLANK(D25),ISBLANK(E25)), """", SUM(C25:E25))"
$ws.Range("H25").Value = 87
$ws.Range("I25").Value = "anger, happiness, neutrality, sadness, surprise, fear"
$ws.Range("J25").Value = "Iranian Persian"

$ws.Range("B26").Select()
